# Update Sheet3's base data table (A20:B36) with the new values for 17-nov.
# Sheet1's CB/CC columns VLOOKUP into Sheet3, which itself VLOOKUPs into
# A20:B36, so updating these 16 cells ripples through both sheets.

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Sheet3")

$ws3.Range("B20").Value = 9.592869920641931
$ws3.Range("B21").Value = 3.1188867291981048
$ws3.Range("B22").Value = 6.9944711270895175
$ws3.Range("B23").Value = 7.614555675745585
$ws3.Range("B24").Value = 6.5866282983148343
$ws3.Range("B26").Value = 11.779434027715201
$ws3.Range("B27").Value = 7.6092121741983911
$ws3.Range("B28").Value = 6.6413192052082621
$ws3.Range("B29").Value = 2.3051702742804281
$ws3.Range("B30").Value = 1.8135000000055927
$ws3.Range("B31").Value = 1.5539369999999233
$ws3.Range("B32").Value = 13.821323680359713
$ws3.Range("B33").Value = 7.8880363053138529
$ws3.Range("B34").Value = 10.606774608120533
$ws3.Range("B35").Value = 7.1079796862285534
$ws3.Range("B36").Value = 43.183440334503089

# Sheet1: add a new "17-nov" snapshot column (CL), one column to the right
# of the existing "16-nov" column (CK), mirroring the recalculated CB/CC
# VLOOKUP results of each row (same pattern as the prior CD..CK columns).

$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("CL1").Value = "17-nov"
$ws1.Range("CL1").NumberFormat = $ws1.Range("CK1").NumberFormat

for ($r = 2; $r -le 18; $r++) {
    $src = $ws1.Cells.Item($r, 80)   # CB column holds the freshly-recalculated value
    $dst = $ws1.Cells.Item($r, 90)   # CL column
    $dst.Value = $src.Value2
    $dst.NumberFormat = $ws1.Cells.Item($r, 89).NumberFormat
}

$ws1.Range("CI11").Select() | Out-Null
